$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.784.09'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '2.462.18'
$ws.Range("E3").Value = '  +0.46%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '''560.14'
$ws.Range("E5").Value = '  -0.69%  '
$ws.Range("D6").Value = '''161.70'
$ws.Range("E6").Value = '  -1.39%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.506'
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("E10").Value = '  +0.53%  '
$ws.Range("D11").Value = '''0.331'
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("D12").Value = '''4.85'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("E13").Value = '  +0.12%  '
$ws.Range("D14").Value = '68.684.22'
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("E15").Value = '  -1.78%  '
$ws.Range("D16").Value = '''23.52'
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("D17").Value = '2.466.84'
$ws.Range("E17").Value = '  -0.39%  '
$ws.Range("D18").Value = '''10.65'
$ws.Range("E18").Value = '  -3.11%  '
$ws.Range("D19").Value = '''334.35'
$ws.Range("E19").Value = '  -2.93%  '
$ws.Range("D20").Value = '''6.90'
$ws.Range("E20").Value = '  -3.44%  '
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("E22").Value = '  +0.03%  '
$ws.Range("E23").Value = '  -0.71%  '
$ws.Range("D24").Value = '''66.56'
$ws.Range("E24").Value = '  -2.27%  '
$ws.Range("D25").Value = '''3.63'
$ws.Range("E25").Value = '  -3.21%  '
$ws.Range("D26").Value = '''8.15'
$ws.Range("E26").Value = '  -1.35%  '
$ws.Range("D27").Value = '0.0₃0812'
$ws.Range("E27").Value = '  -3.53%  '
$ws.Range("D28").Value = '''7.16'
$ws.Range("E28").Value = '  -2.01%  '
$ws.Range("D30").Value = '''428.74'
$ws.Range("E30").Value = '  -2.12%  '
$ws.Range("D31").Value = '''1.13'
$ws.Range("E31").Value = '  -4.19%  '
$ws.Range("E32").Value = '  -4.56%  '
$ws.Range("D33").Value = '''158.86'
$ws.Range("E33").Value = '  +1.44%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -0.62%  '
$ws.Range("D37").Value = '''17.72'
$ws.Range("E37").Value = '  -1.00%  '
$ws.Range("E38").Value = '  -2.43%  '
$ws.Range("D39").Value = '''4.41'
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("E40").Value = '  -4.59%  '
$ws.Range("D41").Value = '''1.07'
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("D42").Value = '''2.05'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").Value = '''3.34'
$ws.Range("E43").Value = '  -0.81%  '
$ws.Range("D44").Value = '''130.26'
$ws.Range("E44").Value = '  -3.72%  '
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '''0.482'
$ws.Range("E46").Value = '  -1.09%  '
$ws.Range("D47").Value = '''0.559'
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("D48").Value = '''0.0907'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("E49").Value = '  +0.20%  '
$ws.Range("E50").Value = '  -3.32%  '
$ws.Range("D51").Value = '''4.96'
$ws.Range("E51").Value = '  -8.29%  '
